# Apply the w11.docx content edits described by the commit diff:
#  1. "1. Think" -> "🧠 Think"
#  2. "2. Read"  -> "📖 Read"
#  3. "Syllabus" list item becomes a hyperlink "Course Syllabus & Website"
#  4. A new callout-tip paragraph (style FirstParagraph) is added after it
#  5. A closing ":::" paragraph (style BodyText) is added after that
#
# Note: emoji characters are surrogate pairs in the .docx character model,
# so a Range object's cached Start/End can drift after an emoji-bearing
# InsertAfter. To stay safe we always re-fetch a *fresh* Range from the
# owning paragraph right before each append instead of reusing/collapsing
# one stale Range object repeatedly.

$d = $word.ActiveDocument

# --- 1 & 2: heading emoji prefixes -----------------------------------
$d.Content.Find.Execute("1. Think", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "🧠 Think", 2) | Out-Null

$d.Content.Find.Execute("2. Read", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "📖 Read", 2) | Out-Null

# --- 3: turn the "Syllabus" list paragraph into a hyperlink -----------
# "Syllabus" is (still) the last paragraph in the document body.
$syllabusPara = $d.Paragraphs($d.Paragraphs.Count)

$linkText = "Course Syllabus & Website"
$start = $syllabusPara.Range.Start
$endNoMark = $syllabusPara.Range.End - 1

$textRng = $d.Range($start, $endNoMark)
$textRng.Text = $linkText

$linkRng = $d.Range($start, $start + $linkText.Length)
$d.Hyperlinks.Add($linkRng, "https://example.com/syllabus", "", "", $linkText) | Out-Null

# Re-fetch the (now hyperlinked) paragraph as the anchor for subsequent
# inserts -- it is still the last paragraph in the body.
$syllabusPara = $d.Paragraphs($d.Paragraphs.Count)

# --- 4: new "::: {.callout-tip}" explanatory paragraph ----------------
$anchorRng = $syllabusPara.Range
$anchorRng.Collapse(0) | Out-Null
$anchorRng.InsertParagraphAfter()

$tipPara = $d.Paragraphs($d.Paragraphs.Count)
$tipPara.Style = "FirstParagraph"

function Append-ToTip($text) {
    $rr = $tipPara.Range
    $rr.Collapse(0) | Out-Null
    $rr.InsertAfter($text)
}

Append-ToTip "::: {.callout-tip}"
Append-ToTip " "
Append-ToTip "### Tip"
Append-ToTip " "
Append-ToTip "-"
Append-ToTip " "
Append-ToTip "“📖 Read”"
Append-ToTip ","
Append-ToTip " "
Append-ToTip "“🎧 Listen”"
Append-ToTip ", and/or"
Append-ToTip " "
Append-ToTip "“📺 Watch”"
Append-ToTip " "
Append-ToTip "items are required content for the day, and should be read/heard/watched before class on that day."
Append-ToTip " "
Append-ToTip "-"
Append-ToTip " "
Append-ToTip "“🌐 Browse”"
Append-ToTip " "
Append-ToTip "items should be briefly looked at but do not need to be read deeply."
Append-ToTip " "
Append-ToTip "-"
Append-ToTip " "
Append-ToTip "“📚 Additional Resources”"
Append-ToTip " "
Append-ToTip "do not need to be looked at; they are there to serve, if useful, as further references for your debates, final projects, and general edification later."

# --- 5: closing ":::" paragraph ---------------------------------------
$closeAnchor = $tipPara.Range
$closeAnchor.Collapse(0) | Out-Null
$closeAnchor.InsertParagraphAfter()

$closePara = $d.Paragraphs($d.Paragraphs.Count)
$closePara.Style = "BodyText"
$closeRng = $closePara.Range
$closeRng.Collapse(0) | Out-Null
$closeRng.InsertAfter(":::")

Write-Output "edit complete"
